# Add a new dataset row (row 5) to the DATASETS sheet for the daily
# macro instruments workbook ingest.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATASETS")

$ws.Range("A5").Value = "DAILY_MACRO_INSTRUMENTS_XLSX"
$ws.Range("B5").Value = "daily_series_wide"
$ws.Range("C5").Value = "xlsx"
$ws.Range("D5").Value = "E:\BacktestData\raw\Macro_Instruments.xlsx"
$ws.Range("E5").Value = "full_refresh"
$ws.Range("G5").Value = "America/New_York"
$ws.Range("H5").Value = "1D"
$ws.Range("I5").Value = "close"
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = "daily_series"
$ws.Range("L5").Value = "year"

$notes = "Macro instruments workbook. Multi-sheet ingest. First column is date. Drop metadata rows where date cell equals 'DATES' or date cannot be parsed. If sheet has 5 columns total (date + 4 data): base=open, .1=high, .2=low, .3=close. If sheet has 6 columns total (date + 5 data): base=open, .1=high, .2=low, .3=close, .4=volume. If sheet has 2 columns total (date + 1 data): treat as last. Drop rows where all mapped data columns are NA (removes weekends/holidays). Do not forward fill. series_id = '<sheet>|<field>'.`nseries_id_prefix_mode: first_data_col`nrequired_series_ids: SPX Index|close, VIX Index|close, SX5E Index|close, USGG10YR Index|close"

$ws.Range("M5").Value = $notes

Write-Host "Row 5 added to DATASETS sheet."
